# "optimizing flexibility, day-per-day LNS, validating solutions"
#
# Adds a "Palettes" column (H) to the pickup-points sheet, fills it in for
# every existing row, tweaks LECLERC ROUFFIAC's pickup days/product type,
# and appends two new pickup points (INTER LA VACHE, METRO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Palettes" column header -----------------------------------------
$ws.Range("H1").Value = "Palettes"

# --- Row 6 (LECLERC ROUFFIAC): extra pickup day + product type change -----
$ws.Range("E6").Value = "Lundi, Mardi, Mercredi, Vendredi"
$ws.Range("G6").Value = "A"

# --- Palettes values for existing rows 2-8 ---------------------------------
$ws.Range("H2").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 2
$ws.Range("H7").Value = 2
$ws.Range("H8").Value = 10

# --- New row 9: INTER LA VACHE ---------------------------------------------
$ws.Range("A2:G2").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Range("A9").Value = "INTER LA VACHE"
$ws.Range("B9").Value = "138 Av. de Fronton, 31200 Toulouse"
$ws.Range("C9").Value = 43.6367727429365
$ws.Range("D9").Value = 1.43413198135102
$ws.Range("E9").Value = "Lundi, Mardi, Mercredi, Jeudi, Vendredi"
$ws.Range("F9").Value = 120
$ws.Range("G9").Value = "F"
$ws.Range("H9").Value = 1
$ws.Rows(9).RowHeight = 19.5

# --- New row 10: METRO -------------------------------------------------
$ws.Range("A2:G2").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A10").Value = "METRO"
$ws.Range("B10").Value = "ZI la Glacière, 8 impasse Camo, 31018 Toulouse"
$ws.Range("C10").Value = 43.642620577339002
$ws.Range("D10").Value = 1.4206436003434
$ws.Range("E10").Value = "Lundi, Mardi, Mercredi, Jeudi, Vendredi"
$ws.Range("F10").Value = 120
$ws.Range("G10").Value = "F"
$ws.Range("H10").Value = 1
$ws.Rows(10).RowHeight = 19.5

$ws.Range("H11").Select()
